$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250, shifting existing rows 250-264 down to 251-265
$ws.Rows.Item(250).Insert()

# Populate the newly inserted row 250 with the new record
$ws.Cells.Item(250, 1).Value = 4
$ws.Cells.Item(250, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(250, 3).Value = "Los Lagos"
$ws.Cells.Item(250, 4).Value = "6/2/2022"
$ws.Cells.Item(250, 5).Value = 10
$ws.Cells.Item(250, 6).Value = 100112021
$ws.Cells.Item(250, 7).Value = "Ají"
$ws.Cells.Item(250, 8).Value = "Inferno"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 60
$ws.Cells.Item(250, 11).Value = 30000
$ws.Cells.Item(250, 12).Value = 30000
$ws.Cells.Item(250, 13).Value = 30000
$ws.Cells.Item(250, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(250, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(250, 16).Value = 2500
$ws.Cells.Item(250, 17).Value = 12
$ws.Cells.Item(250, 18).Value = "Hortaliza"
